# Append the latest daily profit data row (2025-12-09) to the sheet,
# mirroring the format of the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

# Column A holds a date-like label that must stay plain text (matching the
# other rows), not get auto-converted into a date serial number.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "12/09/2025"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value = 13162.25
$ws.Cells.Item($row, 3).Value = 0.1905027776609068
$ws.Cells.Item($row, 4).Value = 0.8094972223390932
$ws.Cells.Item($row, 5).Value = -102.51
$ws.Cells.Item($row, 6).Value = -23.57
$ws.Cells.Item($row, 7).Value = -19836.95
$ws.Cells.Item($row, 8).Value = -65.06
$ws.Cells.Item($row, 9).Value = -395.4
$ws.Cells.Item($row, 10).Value = -13.62
